$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "CHUNGA DE LA CRUZ ROSA LILIANA"
$ws.Range("A5").Value = "RUIZ CARRASCO HILLARY SAMANTHA"

$ws.Range("A10").Value = "PALMA CARMENES DE MENA MERCEDES EVERJISTA"
$ws.Range("A11").Value = "GIRON SILUPU JUAN FRANCISCO"

$ws.Range("B2").Value = 194
$ws.Range("B3").Value = 133
$ws.Range("B4").Value = 132
$ws.Range("B5").Value = 131
$ws.Range("B6").Value = 128
$ws.Range("B7").Value = 109
$ws.Range("B8").Value = 91
$ws.Range("B9").Value = 87
$ws.Range("B10").Value = 82
$ws.Range("B11").Value = 79
$ws.Range("B12").Value = 78
$ws.Range("B13").Value = 64
